$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they stay text (matches source inlineStr data),
# mirroring how the original feed writes plain numeric strings without locale/number coercion.
$ws.Range("D2").Value = '29.339.33'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.869.77'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7149'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.06'
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07899'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3087'
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.34'
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08251'
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7232'
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.243'
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.853.22'
$ws.Range("E14").Value = '  -7.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.74'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '29.368.38'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.841'
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '243.83'
$ws.Range("E18").Value = '  +2.02%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '2.117.13'
$ws.Range("E21").Value = '  -5.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.997'
$ws.Range("E23").Value = '  +5.48%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1599'
$ws.Range("E25").Value = '  +11.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.55'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.983'
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.349'
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.377'
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.100'
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05188'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.940'
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.187'
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7232'
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01855'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.695'
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("D40").Value = '1.172.37'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9047'
$ws.Range("E41").Value = '  -2.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.129'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.54'
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.96'
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.016.22'
$ws.Range("E46").Value = '  -7.41%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5288'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.786'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("E49").Value = '  +5.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.260'
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4278'
$ws.Range("E51").Value = '  -0.14%  '
